$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Transporte "
$ws.Range("B8").Value = "R$ 49.780"

$ws.Range("B9").Select()
